$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 1883
$ws.Range("I3").Value = 1991
$ws.Range("I4").Value = 501
$ws.Range("I5").Value = 175
$ws.Range("H6").Value = 7922
$ws.Range("I6").Value = 2417
$ws.Range("H7").Value = 25970
$ws.Range("I7").Value = 6967

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 69
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 125

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 94
$ws.Range("I7").Value = 268

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 43
$ws.Range("I3").Value = 47
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I4").Value = 33
$ws.Range("I8").Value = 429
$ws.Range("I9").Value = 40
$ws.Range("I10").Value = 55
$ws.Range("I11").Value = 122
$ws.Range("I15").Value = 88
$ws.Range("I19").Value = 204
$ws.Range("I21").Value = 50
$ws.Range("I22").Value = 19
$ws.Range("I27").Value = 64
$ws.Range("I29").Value = 452
$ws.Range("I31").Value = 68
$ws.Range("I33").Value = 327
$ws.Range("I37").Value = 223
$ws.Range("I41").Value = 31
$ws.Range("I42").Value = 232
$ws.Range("I43").Value = 64
$ws.Range("H48").Value = 334
$ws.Range("I48").Value = 72
$ws.Range("I54").Value = 154
$ws.Range("I57").Value = 21
$ws.Range("I63").Value = 27
$ws.Range("I64").Value = 72
$ws.Range("I65").Value = 163
$ws.Range("I67").Value = 268
$ws.Range("I75").Value = 27
$ws.Range("I76").Value = 113
$ws.Range("I79").Value = 176
$ws.Range("I83").Value = 131
$ws.Range("I85").Value = 331
$ws.Range("I88").Value = 52
$ws.Range("I90").Value = 83
$ws.Range("I91").Value = 80
$ws.Range("I94").Value = 55
$ws.Range("I95").Value = 116
$ws.Range("I98").Value = 46
$ws.Range("I99").Value = 125
$ws.Range("H101").Value = 25970
$ws.Range("I101").Value = 6967

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 46
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 81
$ws.Range("I7").Value = 327

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 34
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 154

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 137
$ws.Range("I3").Value = 146
$ws.Range("I7").Value = 452

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 85
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 9
$ws.Range("H6").Value = 182
$ws.Range("I6").Value = 38
$ws.Range("H7").Value = 334
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 81
$ws.Range("I3").Value = 128
$ws.Range("I7").Value = 331

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 31

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 80
$ws.Range("I4").Value = 22
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 232

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 27
$ws.Range("I6").Value = 22
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 50
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 26
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 88

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 53
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I6").Value = 12
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 128
$ws.Range("I5").Value = 14
$ws.Range("I6").Value = 145
$ws.Range("I7").Value = 429

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 27

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 83

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I6").Value = 3
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 11
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 19

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 33
